$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: set cell values (row-major order so new shared strings land at the expected indices) ---
$ws.Range("A56").Value = 5
$ws.Range("B56").Value = "Delay"
$ws.Range("C56").Value = 2
$ws.Range("A57").Value = 5
$ws.Range("B57").Value = "PlaySound"
$ws.Range("C57").Value = "Sound"
$ws.Range("D57").Value = "home"
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Talk"
$ws.Range("C58").Value = "Right"
$ws.Range("D58").Value = "少女"
$ws.Range("E58").Value = "你現在要出發去收集維修飛船的零件了，是吧？"
$ws.Range("F58").Value = "WaitInput"
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = "PlaySound"
$ws.Range("C59").Value = "Act"
$ws.Range("D59").Value = "city01"
$ws.Range("A60").Value = 5
$ws.Range("B60").Value = "Talk"
$ws.Range("C60").Value = "Left"
$ws.Range("D60").Value = "少女"
$ws.Range("E60").Value = "就算我不是Lambda-42系列，我也看得出來，這艘船在航入太陽系前就會自動解體"
$ws.Range("F60").Value = "WaitInput"
$ws.Range("A61").Value = 5
$ws.Range("B61").Value = "PlaySound"
$ws.Range("C61").Value = "Act"
$ws.Range("D61").Value = "city02"
$ws.Range("A62").Value = 5
$ws.Range("B62").Value = "Talk"
$ws.Range("C62").Value = "Left"
$ws.Range("D62").Value = "少女"
$ws.Range("E62").Value = "你還在等什麼呢？收整好裝備後就出發吧"
$ws.Range("F62").Value = "WaitInput"
$ws.Range("A64").Value = 6
$ws.Range("B64").Value = "PlaySound"
$ws.Range("C64").Value = "Sound"
$ws.Range("D64").Value = "home"
$ws.Range("A65").Value = 6
$ws.Range("B65").Value = "Talk"
$ws.Range("C65").Value = "Right"
$ws.Range("D65").Value = "少女"
$ws.Range("E65").Value = "你應該也注意到了，我暫時改裝了飛船的能量引擎"
$ws.Range("F65").Value = "WaitInput"
$ws.Range("A66").Value = 6
$ws.Range("B66").Value = "PlaySound"
$ws.Range("C66").Value = "Act"
$ws.Range("D66").Value = "city01"
$ws.Range("A67").Value = 6
$ws.Range("B67").Value = "Talk"
$ws.Range("C67").Value = "Left"
$ws.Range("D67").Value = "少女"
$ws.Range("E67").Value = "收集夠多的素材，我就能使用能量引擎升級你的性能，協助你打倒外頭的敵人"
$ws.Range("F67").Value = "WaitInput"
$ws.Range("A68").Value = 6
$ws.Range("B68").Value = "PlaySound"
$ws.Range("C68").Value = "Act"
$ws.Range("D68").Value = "city02"
$ws.Range("A69").Value = 6
$ws.Range("B69").Value = "Talk"
$ws.Range("C69").Value = "Left"
$ws.Range("D69").Value = "少女"
$ws.Range("E69").Value = "功能越強，完成任務的機率越高，所有的機器人都是這樣"
$ws.Range("F69").Value = "WaitInput"
$ws.Range("A71").Value = 7
$ws.Range("B71").Value = "PlaySound"
$ws.Range("C71").Value = "Sound"
$ws.Range("D71").Value = "home"
$ws.Range("A72").Value = 7
$ws.Range("B72").Value = "Talk"
$ws.Range("C72").Value = "Right"
$ws.Range("D72").Value = "少女"
$ws.Range("E72").Value = "首先前往溫室吧，那裏可能有備份的氧氣循環機"
$ws.Range("F72").Value = "WaitInput"
$ws.Range("A73").Value = 7
$ws.Range("B73").Value = "PlaySound"
$ws.Range("C73").Value = "Act"
$ws.Range("D73").Value = "city01"
$ws.Range("A74").Value = 7
$ws.Range("B74").Value = "Talk"
$ws.Range("C74").Value = "Left"
$ws.Range("D74").Value = "少女"
$ws.Range("E74").Value = "原則上我們不需要氧氣，但這部飛船原意是設計給人類搭乘，所以氧氣循環機是重要的設施"
$ws.Range("F74").Value = "WaitInput"
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "PlaySound"
$ws.Range("C75").Value = "Act"
$ws.Range("D75").Value = "city02"
$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Talk"
$ws.Range("C76").Value = "Left"
$ws.Range("D76").Value = "少女"
$ws.Range("E76").Value = "可惜，要是這份飛船的設計再好一點，就可以忽略這份構造，不必修復它"
$ws.Range("F76").Value = "WaitInput"
$ws.Range("A77").Value = 7
$ws.Range("B77").Value = "PlaySound"
$ws.Range("C77").Value = "Act"
$ws.Range("D77").Value = "city02"
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Talk"
$ws.Range("C78").Value = "Left"
$ws.Range("D78").Value = "少女"
$ws.Range("E78").Value = "設計師或許從來沒想過搭乘飛船的「人」不需要氧氣吧，呵呵"
$ws.Range("F78").Value = "WaitInput"

# --- Step 2: apply styles by pasting formats from cells that already carry the target style ---
$ws.Range("A1").Copy()
foreach ($addr in @("D56", "E56", "F56", "B57", "C57", "D57", "E57", "F57", "B58", "F58", "B59", "C59", "D59", "E59", "F59", "B60", "C60", "F60", "B61", "C61", "D61", "E61", "F61", "B62", "C62", "F62", "B63", "C63", "D63", "E63", "F63", "B64", "C64", "D64", "E64", "F64", "B65", "F65", "B66", "C66", "D66", "E66", "F66", "B67", "C67", "F67", "B68", "C68", "D68", "E68", "F68", "B69", "C69", "F69", "B70", "C70", "D70", "E70", "F70", "B71", "C71", "D71", "E71", "F71", "B72", "F72", "B73", "C73", "D73", "E73", "F73", "B74", "C74", "F74", "B75", "C75", "D75", "E75", "F75", "B76", "C76", "F76", "B77", "C77", "D77", "E77", "F77", "B78", "C78", "F78")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("A2").Copy()
foreach ($addr in @("A56", "A57", "A58", "A59", "A60", "A61", "A62", "E62", "A64", "A65", "A66", "A67", "E67", "A68", "A69", "E69", "A71", "A72", "E72", "A73", "A74", "E74", "A75", "A76", "E76", "A77", "A78")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("C8").Copy()
foreach ($addr in @("B56", "C56", "C58", "D58", "E58", "D60", "E60", "D62", "C65", "D65", "D67", "D69", "C72", "D72", "D74", "D76", "D78")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("E50").Copy()
foreach ($addr in @("E79", "E80", "E81", "E82")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("E50").Copy()
foreach ($addr in @("F79", "G79", "H79", "I79", "F80", "G80", "H80", "I80", "F81", "G81", "H81", "I81", "F82", "G82", "H82", "I82")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# --- Step 3: special bordered styles (new in this revision) ---
# E65: top border (thin, black) on top of the base style used by row D8 (fontId=2 / style index 2)
$ws.Range("D8").Copy()
$ws.Range("E65").PasteSpecial(-4122)
$b65 = $ws.Range("E65").Borders.Item(8)
$b65.Color = 0
$b65.LineStyle = 1

# E78: bottom border (thin, black)
$ws.Range("D8").Copy()
$ws.Range("E78").PasteSpecial(-4122)
$b78 = $ws.Range("E78").Borders.Item(9)
$b78.Color = 0
$b78.LineStyle = 1

$excel.CutCopyMode = $false